# Reformat three question cells so the wrapped text fits better on screen:
# insert a manual line break in the question text and turn on WrapText
# for those cells, matching the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6, column B: "For which of these concepts does Python have built in support?"
$ws.Range("B6").Value = "For which of these concepts does Python`nhave built in support?"
$ws.Range("B6").WrapText = $true

# Row 11, column B: "What is the name of the output of a compiler after it translates the program?"
$ws.Range("B11").Value = "What is the name of the output of a compiler`nafter it translates the program?"
$ws.Range("B11").WrapText = $true

# Row 16, column B: "Which of the following commands can be used to bring in code from other files?"
$ws.Range("B16").Value = "Which of the following commands can be used to`n bring in code from other files?"
$ws.Range("B16").WrapText = $true

# Move the active selection to where the author last left it.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
[void]$ws.Range("B18").Select()
